$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header style (bold font, border, centered alignment)
# from H1 onto the two new header cells I1 and J1 so they reuse the same
# cell style ("s=1") instead of creating a brand new style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data columns I and J
$ws.Range("I2").Value = 4
$ws.Range("J2").Value = 6

$ws.Range("I3").Value = 8
$ws.Range("J3").Value = 9

$ws.Range("I4").Value = 8
$ws.Range("J4").Value = 9

$ws.Range("I5").Value = 5
$ws.Range("J5").Value = 6
